$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new player/score row after the existing data (row 39 -> row 40)
$ws.Range("A40").Value = "UGNE"
$ws.Range("B40").Value = 200
